$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing G12 value (2.185 -> 2.6185)
$ws.Range("G12").Value = 2.6185

# Add new data for row 12 (M,N,O,P)
$ws.Range("M12").Value = "K"
$ws.Range("N12").Value = "Stiffness"
$ws.Range("O12").Value = 0.0254
$ws.Range("P12").Value = "(N m)/rad"

# Add new row 13 (M,N,O,P)
$ws.Range("M13").Value = "Kt"
$ws.Range("N13").Value = "Thrust_current_Torque constant"
$ws.Range("O13").Value = 0.0064
$ws.Range("P13").Value = "(N m)/A"

# Update selection on the sheet view to O12
$ws.Range("O12").Select()
